$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes existing rows 3..93 down to 4..94)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new record's data
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44882
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 300000001
$ws.Range("G3").Value = "Rabanito"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 7000
$ws.Range("N3").Value = "`$/docena de paquetes"
$ws.Range("O3").Value = "Provincia de Cautín"
$ws.Range("P3").Value = 583
$ws.Range("Q3").Value = 12
$ws.Range("R3").Value = "Hortaliza"
